$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "NA" for the duplicate_image_filename column (E) for the
# practice + main stimuli rows (rows 2 through 21).
$ws.Range("E2:E21").Value = "NA"
